# Refresh crypto price ("Price") and volume-change ("Volume(1h)") columns
# with the latest scraped values, as produced by the scheduled
# "Updated cryptos list ... with GitHub Actions" workflow run.
#
# The sheet stores these numeric-looking values as literal text (e.g.
# "42.350.47", "40.00", "  +0.53%  ") rather than real numbers, so each
# write forces the cell to Text format first and restores the default
# "Normal" style afterwards to avoid leaving stray number formatting
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $NewValue) {
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}


Set-TextValue "D2" "42.600.37"
Set-TextValue "E2" "  +1.02%  "
Set-TextValue "D3" "2.300.73"
Set-TextValue "E3" "  -0.55%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "316.79"
Set-TextValue "E5" "  +1.41%  "
Set-TextValue "D6" "103.96"
Set-TextValue "E6" "  -1.68%  "
Set-TextValue "E7" "  +0.72%  "
Set-TextValue "E8" "  +0.09%  "
Set-TextValue "D9" "0.612"
Set-TextValue "E9" "  +0.49%  "
Set-TextValue "D10" "40.11"
Set-TextValue "E10" "  -0.13%  "
Set-TextValue "D11" "0.0911"
Set-TextValue "E11" "  -0.29%  "
Set-TextValue "E12" "  +0.75%  "
Set-TextValue "D13" "0.106"
Set-TextValue "E13" "  +0.32%  "
Set-TextValue "D14" "0.967"
Set-TextValue "E15" "  -1.93%  "
Set-TextValue "D16" "2.647.63"
Set-TextValue "E16" "  -0.36%  "
Set-TextValue "D17" "2.295.86"
Set-TextValue "E17" "  -1.01%  "
Set-TextValue "D18" "42.369.63"
Set-TextValue "E18" "  +0.59%  "
Set-TextValue "E20" "  +1.09%  "
Set-TextValue "D21" "73.16"
Set-TextValue "E21" "  -1.90%  "
Set-TextValue "D22" "3.57"
Set-TextValue "E22" "  +2.80%  "
Set-TextValue "D23" "276.92"
Set-TextValue "E23" "  +6.94%  "
Set-TextValue "E24" "  +20.51%  "
Set-TextValue "E25" "  -0.83%  "
Set-TextValue "E26" "  -0.24%  "
Set-TextValue "E27" "  -1.33%  "
Set-TextValue "E28" "  +3.47%  "
Set-TextValue "D29" "22.79"
Set-TextValue "E29" "  +0.13%  "
Set-TextValue "D30" "35.82"
Set-TextValue "E30" "  +0.76%  "
Set-TextValue "D31" "165.28"
Set-TextValue "E31" "  +1.48%  "
Set-TextValue "D32" "0.0875"
Set-TextValue "E32" "  -2.35%  "
Set-TextValue "D33" "5.89"
Set-TextValue "E33" "  +0.67%  "
Set-TextValue "E34" "  +5.53%  "
Set-TextValue "E35" "  -10.66%  "
Set-TextValue "D36" "0.116"
Set-TextValue "E36" "  -2.16%  "
Set-TextValue "E37" "  +5.17%  "
Set-TextValue "E38" "  +1.52%  "
Set-TextValue "E39" "  +3.83%  "
Set-TextValue "D40" "2.76"
Set-TextValue "E40" "  -0.67%  "
Set-TextValue "D42" "69.80"
Set-TextValue "E42" "  -0.94%  "
Set-TextValue "D43" "95.11"
Set-TextValue "E43" "  -3.17%  "
Set-TextValue "E44" "  -0.93%  "
Set-TextValue "E45" "  -0.24%  "
Set-TextValue "D46" "82.29"
Set-TextValue "E46" "  +10.12%  "
Set-TextValue "D47" "12.05"
Set-TextValue "E47" "  -0.63%  "
Set-TextValue "E48" "  +1.20%  "
Set-TextValue "D49" "8.89"
Set-TextValue "E49" "  -1.17%  "
Set-TextValue "D50" "1.589.35"
Set-TextValue "E50" "  +2.58%  "
Set-TextValue "D51" "5.09"
Set-TextValue "E51" "  -5.14%  "
